# Mark the zh-cn / de-de handoff rows as a failed transform:
#  - Status              -> "Handoff transform failed"
#  - Latest Handoff File  -> cleared (hyperlink removed)
#  - Latest Handoff Datetime -> reset to 0001-01-01 00:00:00
#  - Handoff Reason       -> "Ignored"
# The "Latest Handback DateTime" (G) cells were already 0001-01-01 00:00:00.

$wb = $excel.ActiveWorkbook

# The Overview sheet's B2 (zh-cn) / C2 (de-de) cells mirror the same
# "Status" text as the per-language sheets below, via the shared string
# table - update them too so every "Ready for handoff" becomes
# "Handoff transform failed".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

$sheetInfo = @{
    "zh-cn" = @{ MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/47ed704c464b00e844de38a0e5e2ed695b1ce6f6/e2e/194bb190-4575-4b59-900e-1c658473f1f8.md"; ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/47ed704c464b00e844de38a0e5e2ed695b1ce6f6/.localization-config" }
    "de-de" = @{ MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/47ed704c464b00e844de38a0e5e2ed695b1ce6f6/e2e/194bb190-4575-4b59-900e-1c658473f1f8.md"; ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/47ed704c464b00e844de38a0e5e2ed695b1ce6f6/.localization-config" }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetInfo[$sheetName]

    # Status -> failed
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff Datetime reset
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason
    $ws.Range("H2").Value = "Ignored"

    # Drop the "Latest Handoff File" cell + its hyperlink. The engine's
    # Hyperlinks.Delete() clears every hyperlink on the sheet (not just the
    # target range), so remove them all, clear C2, then restore the two
    # hyperlinks (A2 -> source md file, A3 -> .localization-config) that
    # must remain.
    $ws.Hyperlinks.Delete()
    $ws.Range("C2").Clear()
    $ws.Hyperlinks.Add($ws.Range("A2"), $info.MdUrl, "", "", "194bb190-4575-4b59-900e-1c658473f1f8.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $info.ConfigUrl, "", "", ".localization-config") | Out-Null
}
